$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($ws, $ref, $val) {
    $ws.Range($ref).Value = "'" + $val
    $ws.Range($ref).Style = "Normal"
}

Set-TextCell $ws "D2" "27.479.47"
Set-TextCell $ws "E2" "  -0.54%  "
Set-TextCell $ws "D3" "1.828.27"
Set-TextCell $ws "E3" "  -1.89%  "
Set-TextCell $ws "D4" "1.006"
Set-TextCell $ws "E4" "  -0.53%  "
Set-TextCell $ws "D5" "331.75"
Set-TextCell $ws "E5" "  -0.44%  "
Set-TextCell $ws "E6" "  -0.46%  "
Set-TextCell $ws "D7" "0.4576"
Set-TextCell $ws "D8" "0.3806"
Set-TextCell $ws "E8" "  -2.83%  "
Set-TextCell $ws "D9" "46.41"
Set-TextCell $ws "E9" "  +1.24%  "
Set-TextCell $ws "D10" "0.07897"
Set-TextCell $ws "E10" "  -1.04%  "
Set-TextCell $ws "D11" "0.9700"
Set-TextCell $ws "E11" "  -3.17%  "
Set-TextCell $ws "D12" "21.04"
Set-TextCell $ws "E12" "  -3.31%  "
Set-TextCell $ws "B13" "Polkadot"
Set-TextCell $ws "C13" "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextCell $ws "D13" "5.885"
Set-TextCell $ws "E13" "  -1.81%  "
Set-TextCell $ws "B14" "WrappedEther"
Set-TextCell $ws "C14" "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextCell $ws "D14" "1.810.64"
Set-TextCell $ws "E14" "  -2.72%  "
Set-TextCell $ws "D15" "7.069"
Set-TextCell $ws "E15" "  -2.37%  "
Set-TextCell $ws "D16" "1.006"
Set-TextCell $ws "E16" "  -0.55%  "
Set-TextCell $ws "D17" "89.68"
Set-TextCell $ws "E17" "  +1.43%  "
Set-TextCell $ws "D18" "0.06623"
Set-TextCell $ws "E18" "  -1.62%  "
Set-TextCell $ws "E19" "  -1.51%  "
Set-TextCell $ws "D20" "17.11"
Set-TextCell $ws "E20" "  +0.11%  "
Set-TextCell $ws "D21" "1.004"
Set-TextCell $ws "E21" "  -0.40%  "
Set-TextCell $ws "D22" "27.466.16"
Set-TextCell $ws "E22" "  -0.58%  "
Set-TextCell $ws "D23" "5.333"
Set-TextCell $ws "E23" "  -2.21%  "
Set-TextCell $ws "E24" "  -0.66%  "
Set-TextCell $ws "D25" "2.309"
Set-TextCell $ws "E25" "  -0.13%  "
Set-TextCell $ws "D26" "2.038.39"
Set-TextCell $ws "E26" "  -2.34%  "
Set-TextCell $ws "D27" "155.70"
Set-TextCell $ws "E27" "  -2.25%  "
Set-TextCell $ws "D28" "19.38"
Set-TextCell $ws "E28" "  -1.99%  "
Set-TextCell $ws "D29" "2.062"
Set-TextCell $ws "E29" "  -3.78%  "
Set-TextCell $ws "D30" "5.295"
Set-TextCell $ws "E30" "  -2.47%  "
Set-TextCell $ws "D31" "118.52"
Set-TextCell $ws "E31" "  -2.65%  "
Set-TextCell $ws "D32" "0.9425"
Set-TextCell $ws "E32" "  -3.78%  "
Set-TextCell $ws "D33" "0.09304"
Set-TextCell $ws "E33" "  -1.87%  "
Set-TextCell $ws "D34" "3.589"
Set-TextCell $ws "E34" "  -0.78%  "
Set-TextCell $ws "D35" "5.253"
Set-TextCell $ws "E35" "  -0.83%  "
Set-TextCell $ws "D36" "1.327"
Set-TextCell $ws "E36" "  -0.41%  "
Set-TextCell $ws "D37" "0.05938"
Set-TextCell $ws "E37" "  -1.91%  "
Set-TextCell $ws "D38" "0.02179"
Set-TextCell $ws "E38" "  -2.14%  "
Set-TextCell $ws "D39" "8.059"
Set-TextCell $ws "E39" "  -3.03%  "
Set-TextCell $ws "D40" "1.145"
Set-TextCell $ws "E40" "  -4.20%  "
Set-TextCell $ws "D41" "0.5771"
Set-TextCell $ws "E41" "  -3.01%  "
Set-TextCell $ws "D42" "0.1827"
Set-TextCell $ws "E42" "  -2.89%  "
Set-TextCell $ws "D43" "9.987"
Set-TextCell $ws "E43" "  -2.76%  "
Set-TextCell $ws "E44" "  +1.29%  "
Set-TextCell $ws "D45" "0.5451"
Set-TextCell $ws "E45" "  -3.33%  "
Set-TextCell $ws "D46" "11.95"
Set-TextCell $ws "E46" "  -2.55%  "
Set-TextCell $ws "E47" "  -2.61%  "
Set-TextCell $ws "D48" "111.00"
Set-TextCell $ws "E48" "  -0.72%  "
Set-TextCell $ws "D49" "0.06612"
Set-TextCell $ws "E49" "  -2.17%  "
Set-TextCell $ws "E50" "  -0.70%  "
Set-TextCell $ws "D51" "1.043"
Set-TextCell $ws "E51" "  -1.17%  "

Write-Host "All cells updated"